# Update the "取得日時" (acquired timestamp) column for rows 2-10 on the
# "ランサーズ" sheet from "2025-11-15 06:24:44" to "2025-11-15 06:32:08".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-15 06:24:44"
$newTimestamp = "2025-11-15 06:32:08"

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
